# Update the "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# row => new value for column F
$updates = @{
    6  = 39
    7  = 59
    8  = 70
    10 = 1168
    11 = 1478
    13 = 364
    15 = 118
    18 = 100
    19 = 261
    20 = 282
    21 = 306
    22 = 1688
    25 = 168
    26 = 639
    28 = 192
    29 = 4027
    32 = 244
    33 = 1035
    34 = 120
    36 = 255
    38 = 148
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
